# Apply updated "dSF" (column F) values for the rows that were repulled /
# recalculated, as described in the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = -5
    12 = 1
    15 = -1
    22 = -2
    23 = -7
    30 = -5
    35 = -4
    40 = 1
    43 = 1
    45 = 3
    48 = 1
    49 = -1
    54 = 5
    57 = -9
    58 = -5
    61 = -5
    62 = 8
    65 = 2
    66 = -6
    69 = -5
    71 = 4
    72 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
